$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @(
    "AVGO",
    "TSLA",
    "JPM",
    "LLY",
    "V",
    "ORCL",
    "BABA",
    "AMD",
    "PLTR",
    "ABBV",
    "NFLX",
    "GS",
    "JPM-PD",
    "JPM-PC",
    "UNH",
    "MS",
    "BAC-PK",
    "GE",
    "BAC-PL",
    "C",
    "NVO",
    "LRCX",
    "SAP",
    "IBM",
    "BAC-PE",
    "BML-PL",
    "RTX",
    "AMAT",
    "BAC-PB",
    "AXP",
    "TMO",
    "CRM",
    "SHEL",
    "KLAC",
    "WFC-PY",
    "DIS",
    "WFC-PL",
    "APH",
    "ISRG",
    "ABT",
    "SAN",
    "SCHW",
    "BX",
    "APP",
    "ANET",
    "SHOP",
    "ACN",
    "BLK",
    "UBER",
    "DHR",
    "BKNG",
    "QCOM",
    "SPGI",
    "INTU",
    "UBS",
    "PDD",
    "BBVA",
    "NOW",
    "COF",
    "BSX",
    "NEM",
    "SONY",
    "PANW",
    "IBKR",
    "ADBE",
    "WFC-PC",
    "VRTX",
    "PGR",
    "CRWD",
    "MELI",
    "KKR",
    "AEM",
    "CVS",
    "BN",
    "SPOT",
    "ADP",
    "CEG",
    "CVNA",
    "ICE",
    "GD",
    "SNPS",
    "HOOD",
    "NKE",
    "MCO",
    "BCS",
    "DASH",
    "NU",
    "FCX",
    "CDNS",
    "B",
    "HWM",
    "NTES",
    "ELV",
    "ORLY",
    "MS-PK",
    "AMT",
    "BAM",
    "MS-PI",
    "CRH",
    "ABNB",
    "TDG",
    "DB",
    "MS-PF",
    "MS-PE",
    "CMI",
    "DELL",
    "APO",
    "INFY",
    "SE",
    "USB-PP",
    "AON",
    "SNOW",
    "GS-PD",
    "RELX",
    "NWG",
    "MRVL",
    "WBD",
    "VRT",
    "WPM",
    "TEL",
    "NET",
    "AJG",
    "FTNT",
    "AZO",
    "MFC",
    "DUK-PA",
    "RACE",
    "RKT",
    "CTA-PB",
    "URI",
    "NXPI",
    "ADSK",
    "COIN",
    "NDAQ",
    "SRE",
    "IDXX",
    "TRI",
    "AU",
    "BIDU",
    "PYPL",
    "VST",
    "F",
    "CMG",
    "CCJ",
    "RBLX",
    "ARGX",
    "MET",
    "EA",
    "SCHW-PD",
    "WDAY",
    "BSBR",
    "CBRE",
    "FNV",
    "EW",
    "GFI",
    "ARES",
    "CRWV",
    "AXON",
    "DDOG",
    "ALNY",
    "ROK",
    "CTA-PA",
    "HEI",
    "MSTR",
    "KGC",
    "MSCI",
    "TTWO",
    "SPG-PJ",
    "EXC",
    "ROP",
    "JD",
    "EBAY",
    "MET-PE",
    "RKLB",
    "MET-PA",
    "EL",
    "CTSH",
    "TCOM",
    "LVS",
    "IQV",
    "PUK",
    "RDDT",
    "CPRT",
    "XYZ",
    "HEI-A",
    "ALC",
    "GRMN",
    "PEG",
    "ASTS",
    "A",
    "PRU",
    "PSA-PH",
    "PAYX",
    "CCI",
    "MDLN",
    "FICO",
    "VEEV",
    "FISV",
    "TEAM",
    "RYAAY",
    "CPNG",
    "STT",
    "CLS",
    "SYM",
    "APO-PA",
    "XYL",
    "UAL",
    "SATS",
    "LYV",
    "ZS",
    "INSM",
    "NTRA",
    "RJF",
    "MDB",
    "EXPE",
    "CHT",
    "BE",
    "ESLT",
    "PSA-PK",
    "ALL-PH",
    "SOFI",
    "ADM",
    "ALL-PB",
    "FOXA",
    "HUM",
    "EXR",
    "FIS",
    "FOX",
    "VRSK",
    "FLUT",
    "FTAI",
    "BNTX",
    "MTD",
    "NRG",
    "SYF",
    "DXCM",
    "TME",
    "NTRS",
    "STLA",
    "CSGP",
    "PHG",
    "IRM",
    "ALAB",
    "PAAS",
    "WIT",
    "EC",
    "BRO",
    "PPL",
    "NMR",
    "EFX",
    "ES",
    "FSLR",
    "STE",
    "DLTR",
    "AWK",
    "OMC",
    "AVB",
    "VLTO",
    "DLR-PK",
    "GFS",
    "RGLD",
    "CW",
    "EQR",
    "BR",
    "LDOS",
    "PSTG",
    "SQM",
    "FLEX",
    "AXIA-PC",
    "ILMN",
    "VRSN",
    "TPG",
    "TROW",
    "WAT",
    "CRDO",
    "NBIS",
    "LULU",
    "OWL",
    "CNC",
    "FUTU",
    "AFRM",
    "DLR-PJ",
    "TLK",
    "CYBR",
    "FWONK",
    "ALB",
    "PSLV",
    "FWONA",
    "CG",
    "RL",
    "CIB",
    "GPN",
    "AS",
    "SSNC",
    "GMAB",
    "TWLO",
    "Q",
    "HL",
    "SBAC",
    "RCI",
    "CHKP",
    "PTC",
    "NTAP",
    "TOST",
    "GIB",
    "RIVN",
    "PODD",
    "TYL",
    "RVMD",
    "BWXT",
    "KTOS",
    "MRNA",
    "HIG-PG",
    "GRAB",
    "IOT",
    "CX",
    "PHYS",
    "U",
    "SMCI",
    "DKS",
    "HPQ",
    "CRCL",
    "FITBI",
    "XPEV",
    "IT",
    "PSNYW",
    "AGI",
    "ALLY",
    "PNR",
    "PINS",
    "ENTG",
    "ESS",
    "SN",
    "WST",
    "HUBS",
    "NWS",
    "IREN",
    "ZG",
    "FN",
    "TRMB",
    "MEDP",
    "JLL",
    "QXO",
    "Z",
    "HOLX",
    "TRU",
    "TTD",
    "INVH",
    "CDE",
    "NLY",
    "HII",
    "MAA",
    "OKTA",
    "COO",
    "SUI",
    "TLN",
    "GFL",
    "WMG",
    "KSPI",
    "ROKU",
    "ROIV",
    "GEN",
    "NWSA",
    "FFIV",
    "H",
    "IONQ",
    "ONON",
    "CRBG",
    "DKNG",
    "AA",
    "AVAV",
    "MGA",
    "KEY-PI",
    "BBIO",
    "GH",
    "MLI",
    "HMY",
    "PFGC",
    "ULS",
    "GDDY",
    "ARCC",
    "ICLR",
    "FIG",
    "RGC",
    "KRMN",
    "RF-PC",
    "W",
    "ASND",
    "CACI",
    "AKAM",
    "JHX",
    "PEN",
    "EVR",
    "CELH",
    "DPZ",
    "EMBJ",
    "SBSW",
    "LOGI",
    "EQH",
    "GWRE",
    "EG",
    "SOLV",
    "NVMI",
    "BILI",
    "RBRK",
    "FIGR",
    "AMH",
    "RVTY",
    "JKHY",
    "RYAN",
    "PSKY",
    "CHWY",
    "UNM",
    "SNAP",
    "JEF",
    "BNT",
    "OKLO",
    "HLI",
    "EQX",
    "IVZ",
    "AGNC",
    "BMNR",
    "GLXY",
    "AG",
    "GMED",
    "AMKR",
    "IAG",
    "DT",
    "ACGLO",
    "ZBRA",
    "DOC",
    "SMMT",
    "RMBS",
    "EPAM",
    "JOBY",
    "CR",
    "CMA",
    "TEM",
    "WYNN",
    "NYT",
    "AEG",
    "EXEL",
    "NTNX",
    "BSY",
    "DOCU",
    "CPT",
    "BXP",
    "LSCC",
    "STN",
    "CEF",
    "MDGL",
    "BAH",
    "MICC",
    "MP",
    "QGEN",
    "WTRG",
    "SARO",
    "CRL",
    "AR",
    "UHAL",
    "MOH",
    "DRS",
    "MANH",
    "FDS",
    "AFG",
    "CART",
    "SEIC",
    "TECH",
    "CAE",
    "YMM",
    "XP",
    "ONTO",
    "PCOR",
    "KLAR",
    "APLD",
    "JAZZ",
    "HBM",
    "CHYM",
    "VNO-PL",
    "UHAL-B",
    "SAIL",
    "SANM",
    "BVN",
    "NGD",
    "VNO-PM",
    "COMP",
    "BIO-B",
    "TTMI",
    "ALV",
    "SOLS",
    "EGO",
    "DDS",
    "BROS",
    "REXR",
    "AAL",
    "AYI",
    "ABVX",
    "NCLH",
    "IESC",
    "UWMC",
    "ARWR",
    "GS-PC",
    "RGEN",
    "DOX",
    "STEP",
    "MORN",
    "AMG",
    "UEC",
    "LUMN",
    "GGAL",
    "QBTS",
    "TTAN",
    "RZB",
    "AGNCM",
    "AGNCN",
    "PEGA",
    "VIPS",
    "OR",
    "AUR",
    "GGB",
    "UGI",
    "PL",
    "LEVI"
)

for ($i = 0; $i -lt $tickers.Count; $i++) {
    $ws.Cells.Item($i + 7, 1).Value = $tickers[$i]
}

$ws.Range("A530:A532").ClearContents()

